$wb = $excel.ActiveWorkbook
$leadSheet = $wb.Worksheets.Item("CreateNewLead")

$new = $wb.Worksheets.Add($null, $leadSheet)
$new.Name = "CreateNewContact"

$new.Range("A1").Value = "FirstName"
$new.Range("B1").Value = "LastName"
$new.Range("A2").Value = "John"
$new.Range("B2").Value = "Smith"
$new.Range("A1:B1").Font.Bold = $true

$leadSheet.Range("A1:B2").Select()

$new.Activate()
$new.Range("C5").Select()
